$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H40").Value = 7799.143
$ws1.Range("I40").Value = 2297
$ws1.Range("K40").Value = 2297
$ws1.Range("M40").Value = -2122
$ws1.Range("H55").Value = 180
$ws1.Range("I55").Value = 150
$ws1.Range("J55").Value = 300
$ws1.Range("K55").Value = 150
$ws1.Range("L55").Value = 300
$ws1.Range("M55").Value = 64
$ws1.Range("N55").Value = -728
$ws1.Range("H62").Value = 5100
$ws1.Range("I62").Value = 4875
$ws1.Range("J62").Value = 6000
$ws1.Range("K62").Value = 4875
$ws1.Range("L62").Value = 6000
$ws1.Range("M62").Value = -4251
$ws1.Range("N62").Value = -7248
$ws1.Range("H65").Value = 5100
$ws1.Range("I65").Value = 4875
$ws1.Range("J65").Value = 6000
$ws1.Range("K65").Value = 24375
$ws1.Range("L65").Value = 30000
$ws1.Range("M65").Value = -21255
$ws1.Range("N65").Value = -36240
$ws1.Range("H97").Value = 8026.5
$ws1.Range("J97").Value = 8026.5
$ws1.Range("L97").Value = 24079.5
$ws1.Range("N97").Value = -25071.5
$ws1.Range("H98").Value = 638.25
$ws1.Range("J98").Value = 498.5
$ws1.Range("L98").Value = 498.5
$ws1.Range("N98").Value = -3494.5
$ws1.Range("H122").Value = 638.25
$ws1.Range("J122").Value = 498.5
$ws1.Range("L122").Value = 1495.5
$ws1.Range("N122").Value = -6395.5
$ws1.Range("H127").Value = 724.7143
$ws1.Range("I127").Value = 724.7143
$ws1.Range("K127").Value = 2174.1429
$ws1.Range("M127").Value = 2785.8571

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H37").Value = 18333.334
$ws2.Range("H63").Value = 3850
$ws2.Range("I63").Value = 3850
$ws2.Range("K63").Value = 3850
$ws2.Range("M63").Value = -3164
$ws2.Range("H66").Value = 3850
$ws2.Range("I66").Value = 3850
$ws2.Range("K66").Value = 19250
$ws2.Range("M66").Value = -15818

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H86").Value = 6399.875
$ws3.Range("I86").Value = 2400
$ws3.Range("K86").Value = 2400
$ws3.Range("M86").Value = -1277
$ws3.Range("H89").Value = 6399.875
$ws3.Range("I89").Value = 2400
$ws3.Range("K89").Value = 12000
$ws3.Range("M89").Value = -6384
$ws3.Range("H105").Value = 2063.625
$ws3.Range("I105").Value = 2130
$ws3.Range("K105").Value = 2130
$ws3.Range("M105").Value = -383
$ws3.Range("H135").Value = 0
$ws3.Range("J135").Value = 0
$ws3.Range("L135").ClearContents()
$ws3.Range("N135").Value = 0

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H50").Value = 32857.145
$ws4.Range("J50").Value = 34000
$ws4.Range("L50").Value = 34000
$ws4.Range("N50").Value = -35250
$ws4.Range("H62").Value = 4750
$ws4.Range("J62").Value = 4750
$ws4.Range("L62").Value = 4750
$ws4.Range("N62").Value = -5998
$ws4.Range("H65").Value = 4750
$ws4.Range("J65").Value = 4750
$ws4.Range("L65").Value = 23750
$ws4.Range("N65").Value = -29990
$ws4.Range("H99").Value = 2344.4443
$ws4.Range("I99").Value = 2344.4443
$ws4.Range("K99").Value = 2344.4443
$ws4.Range("M99").Value = -846.4443000000001
$ws4.Range("H105").Value = 1262.7059
$ws4.Range("I105").Value = 1211.1333
$ws4.Range("J105").Value = 1649.5
$ws4.Range("K105").Value = 1211.1333
$ws4.Range("L105").Value = 1649.5
$ws4.Range("M105").Value = 535.8667
$ws4.Range("N105").Value = -5143.5
$ws4.Range("H122").Value = 1920
$ws4.Range("I122").Value = 1920
$ws4.Range("K122").Value = 5760
$ws4.Range("M122").Value = -3310
$ws4.Range("H126").Value = 2344.4443
$ws4.Range("I126").Value = 2344.4443
$ws4.Range("K126").Value = 7033.3329
$ws4.Range("M126").Value = -4563.3329

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H12").Value = 37.666668
$ws5.Range("J12").Value = 12.333333
$ws5.Range("L12").Value = 36.999999
$ws5.Range("N12").Value = -382.999999

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H80").Value = 3249.75
$ws6.Range("I80").Value = 2499
$ws6.Range("K80").Value = 2499
$ws6.Range("M80").Value = -1501
$ws6.Range("H83").Value = 3249.75
$ws6.Range("I83").Value = 2499
$ws6.Range("K83").Value = 12495
$ws6.Range("M83").Value = -7503
$ws6.Range("H122").Value = 8342.412
$ws6.Range("I122").Value = 3131.4666
$ws6.Range("J122").Value = 47424.5
$ws6.Range("K122").Value = 9394.399800000001
$ws6.Range("L122").Value = 142273.5
$ws6.Range("M122").Value = -6944.399800000001
$ws6.Range("N122").Value = -147173.5

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H3").Value = 34998
$ws7.Range("I3").Value = 0
$ws7.Range("J3").Value = 34998
$ws7.Range("K3").Value = 0
$ws7.Range("L3").ClearContents()
$ws7.Range("M3").Value = 34998
$ws7.Range("N3").Value = -35222
$ws7.Range("H15").Value = 34998
$ws7.Range("I15").Value = 0
$ws7.Range("J15").Value = 34998
$ws7.Range("K15").Value = 0
$ws7.Range("L15").ClearContents()
$ws7.Range("M15").Value = 34998
$ws7.Range("N15").Value = -35338
$ws7.Range("H16").Value = 3349
$ws7.Range("I16").Value = 3349
$ws7.Range("K16").Value = 3349
$ws7.Range("M16").Value = -3179
$ws7.Range("H93").Value = 5000
$ws7.Range("I93").Value = 5000
$ws7.Range("K93").Value = 5000
$ws7.Range("M93").Value = -3752
$ws7.Range("H132").Value = 22776.908
$ws7.Range("I132").Value = 20054.6
$ws7.Range("J132").Value = 50000
$ws7.Range("K132").Value = 60163.8
$ws7.Range("L132").Value = 150000
$ws7.Range("M132").Value = -57633.8
$ws7.Range("N132").Value = -155060
$ws7.Range("H136").Value = 5549.7144
$ws7.Range("I136").Value = 5369.8
$ws7.Range("K136").Value = 16109.4
$ws7.Range("M136").Value = -13559.4

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H14").Value = 4568
$ws8.Range("I14").Value = 4568
$ws8.Range("K14").Value = 4568
$ws8.Range("M14").Value = -4400
$ws8.Range("H81").Value = 0
$ws8.Range("I81").Value = 0
$ws8.Range("K81").Value = 0
$ws8.Range("M81").ClearContents()
$ws8.Range("H84").Value = 0
$ws8.Range("I84").Value = 0
$ws8.Range("K84").Value = 0
$ws8.Range("M84").ClearContents()
$ws8.Range("H122").Value = 223677.33
$ws8.Range("I122").Value = 286584.56
$ws8.Range("J122").Value = 3502
$ws8.Range("K122").Value = 859753.6799999999
$ws8.Range("L122").Value = 10506
$ws8.Range("M122").Value = -857303.6799999999
$ws8.Range("N122").Value = -15406
$ws8.Range("H132").Value = 3678.68
$ws8.Range("I132").Value = 2998.5
$ws8.Range("K132").Value = 8995.5
$ws8.Range("M132").Value = -6465.5
$ws8.Range("H136").Value = 4001.3333
$ws8.Range("I136").Value = 4002
$ws8.Range("J136").Value = 4000
$ws8.Range("K136").Value = 12006
$ws8.Range("L136").Value = 12000
$ws8.Range("M136").Value = -9456
$ws8.Range("N136").Value = -17100
